# Insert a new weekly price record as row 151 in the daily-logic subset
# sheet, pushing the existing rows 151-183 down to 152-184.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 151; everything at/after 151 shifts down.
$ws.Rows.Item(151).Insert()

# Populate the newly inserted row 151 with the new record.
$ws.Range("A151").Value = 11
$ws.Range("B151").Value = "Vega Monumental Concepción"
$ws.Range("C151").Value = "Bíobío"
$ws.Range("D151").Value = 44943
$ws.Range("D151").NumberFormat = $ws.Range("D152").NumberFormat
$ws.Range("E151").Value = 8
$ws.Range("F151").Value = 100112032
$ws.Range("G151").Value = "Zapallo italiano"
$ws.Range("H151").Value = "Sin especificar"
$ws.Range("I151").Value = "Primera"
$ws.Range("J151").Value = 450
$ws.Range("K151").Value = 6000
$ws.Range("L151").Value = 7000
$ws.Range("M151").Value = 6444
$ws.Range("N151").Value = "$/caja 50 unidades"
$ws.Range("O151").Value = "Región de O'Higgins"
$ws.Range("P151").Value = 129
$ws.Range("Q151").Value = 50
$ws.Range("R151").Value = "Hortaliza"
